# Applies:
#  1. Removes the stray <w:rFonts w:hint="eastAsia"/> from the first
#     paragraph's paragraph-mark run properties (w:pPr/w:rPr).
#  2. Appends two empty paragraphs and a third paragraph containing the
#     text "Sdsdsdksjdksjds dfsdfsdfsfsf" at the end of the document,
#     all tagged with lang="en-US" (matching the document's existing
#     paragraph-mark formatting), with no east-Asia font hint.

$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- Step 1: rewrite the first paragraph, dropping the rFonts hint on
#     the paragraph mark (w:pPr/w:rPr) while leaving everything else,
#     including the run-level rFonts hint on "s", untouched. ---
$firstPara = $d.Paragraphs.First
$firstParaXml = '<w:p xmlns:w="' + $wNs + '" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="16E0BCA0" w14:textId="3E596D7D" w:rsidR="0029050A" w:rsidRPr="00FE2862" w:rsidRDefault="00FE2862">' + `
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Canada has always been known as one of the safest countries in the world, yet when I began my university study in Toronto in 2017, I have heard more crime news since then. By the end of 2020, I moved a new apartment which was located near the Toronto policy headquarter, and ever since then, I heard more and more polic</w:t></w:r>' + `
  '<w:r w:rsidR="008A2207"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>e</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US"/></w:rPr><w:t>s</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">irens. At that point, I started wondering if Toronto was safe to live and </w:t></w:r>' + `
  '<w:r w:rsidR="008A2207"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">if the neighborhood I lived in was secure enough. Thus, I started working on this project and trying to learn more about the crime in Toronto by analyzing the data provided by the Toronto police from 2014 to 2019. </w:t></w:r>' + `
  '</w:p>'
$firstPara.Range.InsertXML($firstParaXml)

# --- Step 2: append the two blank paragraphs + the new text paragraph
#     at the very end of the document body. ---
$endRange = $d.Content
$endRange.Collapse(0)

$pBlank = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
$pText = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Sdsdsdksjdksjds dfsdfsdfsfsf</w:t></w:r></w:p>'

$endRange.InsertXML($pBlank + $pBlank + $pText)
